{"js": "// Replace the worksheet date and every \"AA\u00d7BB=\" multiplication prompt with\n// its new value. Every \"before\" string in this table is unique within the\n// document, so a plain text search + replace is unambiguous for each pair.\nconst replacements = [\n  [\"2026-01-09 Friday\", \"2026-01-10 Saturday\"],\n  [\"68\u00d764=\", \"19\u00d777=\"],\n  [\"47\u00d713=\", \"67\u00d772=\"],\n  [\"41\u00d733=\", \"94\u00d724=\"],\n  [\"53\u00d756=\", \"95\u00d784=\"],\n  [\"43\u00d799=\", \"83\u00d794=\"],\n  [\"76\u00d799=\", \"13\u00d783=\"],\n  [\"58\u00d782=\", \"92\u00d795=\"],\n  [\"47\u00d756=\", \"12\u00d792=\"],\n  [\"93\u00d751=\", \"22\u00d787=\"],\n  [\"31\u00d761=\", \"25\u00d759=\"],\n  [\"37\u00d757=\", \"50\u00d767=\"],\n  [\"80\u00d748=\", \"44\u00d799=\"],\n  [\"32\u00d714=\", \"64\u00d726=\"],\n  [\"38\u00d789=\", \"81\u00d750=\"],\n  [\"74\u00d732=\", \"27\u00d780=\"],\n  [\"95\u00d739=\", \"13\u00d778=\"],\n  [\"36\u00d778=\", \"46\u00d759=\"],\n  [\"66\u00d744=\", \"55\u00d720=\"],\n  [\"36\u00d759=\", \"52\u00d730=\"],\n  [\"17\u00d776=\", \"13\u00d765=\"],\n  [\"83\u00d720=\", \"70\u00d745=\"],\n  [\"48\u00d760=\", \"56\u00d749=\"],\n  [\"54\u00d737=\", \"20\u00d768=\"],\n  [\"24\u00d753=\", \"74\u00d715=\"],\n  [\"54\u00d773=\", \"34\u00d766=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the worksheet date and every \"AA\u00d7BB=\" multiplication prompt with\n# its new value. Every \"before\" string in this table is unique within the\n# document, so Find/Replace on the whole document story is unambiguous for\n# each pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2026-01-09 Friday\", \"2026-01-10 Saturday\"),\n  @(\"68\u00d764=\", \"19\u00d777=\"),\n  @(\"47\u00d713=\", \"67\u00d772=\"),\n  @(\"41\u00d733=\", \"94\u00d724=\"),\n  @(\"53\u00d756=\", \"95\u00d784=\"),\n  @(\"43\u00d799=\", \"83\u00d794=\"),\n  @(\"76\u00d799=\", \"13\u00d783=\"),\n  @(\"58\u00d782=\", \"92\u00d795=\"),\n  @(\"47\u00d756=\", \"12\u00d792=\"),\n  @(\"93\u00d751=\", \"22\u00d787=\"),\n  @(\"31\u00d761=\", \"25\u00d759=\"),\n  @(\"37\u00d757=\", \"50\u00d767=\"),\n  @(\"80\u00d748=\", \"44\u00d799=\"),\n  @(\"32\u00d714=\", \"64\u00d726=\"),\n  @(\"38\u00d789=\", \"81\u00d750=\"),\n  @(\"74\u00d732=\", \"27\u00d780=\"),\n  @(\"95\u00d739=\", \"13\u00d778=\"),\n  @(\"36\u00d778=\", \"46\u00d759=\"),\n  @(\"66\u00d744=\", \"55\u00d720=\"),\n  @(\"36\u00d759=\", \"52\u00d730=\"),\n  @(\"17\u00d776=\", \"13\u00d765=\"),\n  @(\"83\u00d720=\", \"70\u00d745=\"),\n  @(\"48\u00d760=\", \"56\u00d749=\"),\n  @(\"54\u00d737=\", \"20\u00d768=\"),\n  @(\"24\u00d753=\", \"74\u00d715=\"),\n  @(\"54\u00d773=\", \"34\u00d766=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    Write-Output \"WARNING: no match found for $oldText\"\n  }\n}\n"}
